# sdc-generic-workflow.pptx — "Update diagram to use proper terms"
#
# The deck's Shape.Left/Top/Width/Height are single-precision (Single)
# points, same as real PowerPoint's COM model, so a naive EMU/12700
# conversion can truncate one EMU low when the value round-trips back
# through XML. Emu2Pt() nudges the point value up by a fraction of an
# EMU (in 0.05-EMU steps) until the float32 round-trip lands back in
# [emu, emu+1), so the saved <a:off>/<a:ext> integers come out exact.
function Emu2Pt {
    param([double]$emu)
    for ($n = 0; $n -lt 60; $n++) {
        $v = ($emu + $n * 0.05) / 12700.0
        $single = [float]$v
        $back = [double]$single * 12700.0
        if ($back -ge $emu -and $back -lt ($emu + 1)) {
            return $v
        }
    }
    return $emu / 12700.0
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape 7 ("6. EHR system allows for Provider data-entry into the
#     form/template") -- resize/reposition + reword tail -------------
$sh = $s.Shapes.Item(7)
$sh.Left   = Emu2Pt 1015155
$sh.Top    = Emu2Pt 5158581
$sh.Width  = Emu2Pt 1334335
$sh.Height = Emu2Pt 676910
$tr = $sh.TextFrame.TextRange
$tr.Characters(26, 42).Text = "data-entry and correction of Questionnaire Response"

# --- Shape 9 ("7. The EHR system transmits completed structured data
#     in standard format") -------------------------------------------
$sh = $s.Shapes.Item(9)
$tr = $sh.TextFrame.TextRange
$tr.Characters(39, 34).Text = "Questionnaire Response"

# --- Shape 10 ("8. The External Data Repository receives the
#     structured data") -----------------------------------------------
$sh = $s.Shapes.Item(10)
$tr = $sh.TextFrame.TextRange
$tr.Characters(46, 15).Text = "Questionnaire Response"

# --- Shape 11 ("1. EHR system sends request for form/template" /
#     (blank) / "OPTIONAL: Sends some patient data") -- only the first
#     paragraph changes ------------------------------------------------
$sh = $s.Shapes.Item(11)
$tr = $sh.TextFrame.TextRange
$tr.Characters(1, 45).Text = "1. EHR system sends request for populated Questionnaire Response"

# --- Shape 12 ("5. EHR system displays correct form/template") -------
$sh = $s.Shapes.Item(12)
$tr = $sh.TextFrame.TextRange
$tr.Characters(24, 21).Text = "partially completed Questionnaire Response"

# --- Shape 13 ("3. Form/Template repository sends correct
#     form/template") --------------------------------------------------
$sh = $s.Shapes.Item(13)
$tr = $sh.TextFrame.TextRange
$tr.Characters(35, 21).Text = "partially populated Questionnaire Response"

# --- Shape 14 ("4(a). CONDITIONAL FUNCTIONALITY" / "Form is
#     auto-populated with some EHR-derived patient data") -- edit the
#     rightmost span first so the earlier offset stays valid ---------
$sh = $s.Shapes.Item(14)
$tr = $sh.TextFrame.TextRange
$tr.Characters(56, 22).Text = "with EHR-derived "
$tr.Characters(33, 5).Text  = "Questionnaire Response "

# --- Shape 15 ("2. Form/Template repository receives request for
#     form/template") -- whole (single-run) paragraph reworded -------
$sh = $s.Shapes.Item(15)
$tr = $sh.TextFrame.TextRange
$tr.Characters(1, 62).Text = "2. Form repository receives request for populated Questionnaire Response"

# --- Shape 16 ("9. The External Data Repository stores the
#     structured data in standard format") ----------------------------
$sh = $s.Shapes.Item(16)
$tr = $sh.TextFrame.TextRange
$tr.Characters(1, 77).Text = "9. The External Data Repository stores the Questionnaire Response in standard format"

# --- Connector 21 (the little arrow right below shape 7) -------------
$sh = $s.Shapes.Item(21)
$sh.Left   = Emu2Pt 1682323
$sh.Top    = Emu2Pt 5094398
$sh.Width  = Emu2Pt 0
$sh.Height = Emu2Pt 64183

# --- Connector 22 (the arrow below that one) --------------------------
$sh = $s.Shapes.Item(22)
$sh.Left   = Emu2Pt 1682323
$sh.Top    = Emu2Pt 5835491
$sh.Width  = Emu2Pt 0
$sh.Height = Emu2Pt 154550

# --- Shape 26 ("4. EHR System receives correct form/template") -------
$sh = $s.Shapes.Item(26)
$tr = $sh.TextFrame.TextRange
$tr.Characters(24, 21).Text = "Questionnaire Response"

# --- Shape 28 ("2(a). CONDITIONAL FUNCTIONALITY" / "Form is
#     pre-populated with some EHR-provided patient data") ------------
$sh = $s.Shapes.Item(28)
$tr = $sh.TextFrame.TextRange
$tr.Characters(55, 10).Text = "with "
$tr.Characters(33, 5).Text  = "Questionnaire Response "
